$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 7 and 8
$ws.Range("D7").Value = 44511
$ws.Range("D8").Value = 44504
$ws.Range("J8").Value = 500

# Insert new row 9 (duplicate of the former row 8 data, before the J/D edits)
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C9").Value = "Ñuble"
$ws.Range("D9").Value = 44505
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = 300000000
$ws.Range("G9").Value = "Espárragos"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 440
$ws.Range("K9").Value = 900
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = 950
$ws.Range("N9").Value = "$/kilo"
$ws.Range("O9").Value = "Provincia de Diguillín"
$ws.Range("P9").Value = 950
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = "Hortaliza"

# Match style (date number format) used in D2:D8 for the new D9 cell
$ws.Range("D9").NumberFormat = $ws.Range("D8").NumberFormat
